$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Edit 1: extend ". Оценки получились завышенными." with a sentence that
# points to the GitHub notebook used for training, split into three runs
# (the URL itself carries no explicit run properties, matching the target
# OOXML) instead of one big run.
# -------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(". Оценки получились завышенными.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: target sentence not found"
}
$target1 = $d.Range($r1.Start, $r1.End)
$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>. Оценки получились завышенными. Посмотреть о деталях обучения (например, как менялись значения весовых коэффициентов) и модели можно на github по ссылке (</w:t>
      </w:r>
      <w:r>
        <w:t>https://github.com/MrKozelberg/nn-based-multistate-solver-for-se/blob/main/src/nnbmss_for_static_se.ipynb)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target1.InsertXML($xml1)

# -------------------------------------------------------------------------
# Edit 2: " и  посмотреть, ... размерности 4. " becomes
# " и " + oMath(w_{b,max}) + " посмотреть, ... размерности 4. "
# mirroring the w_{a,max} equation that already sits earlier in the same
# paragraph.
# -------------------------------------------------------------------------
$r2 = $d.Content
$searchText2 = " и  посмотреть, как модель работает с большим числом состояний. После этого, перейду к размерности 4. "
$found2 = $r2.Find.Execute($searchText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: target text not found"
}
$target2 = $d.Range($r2.Start, $r2.End)
$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> и </w:t>
      </w:r>
      <m:oMathPara>
        <m:oMath>
          <m:sSub>
            <m:e>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="24"/>
                </w:rPr>
                <m:t>w</m:t>
              </m:r>
            </m:e>
            <m:sub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="24"/>
                </w:rPr>
                <m:t>b,</m:t>
              </m:r>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="24"/>
                </w:rPr>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>max</m:t>
              </m:r>
            </m:sub>
          </m:sSub>
        </m:oMath>
      </m:oMathPara>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> посмотреть, как модель работает с большим числом состояний. После этого, перейду к размерности 4. </w:t>
      </w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target2.InsertXML($xml2)

Write-Host "Applied both edits."
